$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing row 2 values
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 0.1509656666666667
$ws.Range("H2").Value = 0.452897
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 0.7885686666666668
$ws.Range("N2").Value = 2.365706
$ws.Range("O2").Value = 0.4566863346753138
$ws.Range("P2").Value = 0.4566863346753137
$ws.Range("Q2").Value = 0.1190467944757778
$ws.Range("R2").Value = 1.071421150282
$ws.Range("S2").Value = 0.4566863346753138
$ws.Range("T2").Value = 0.4566863346753137

# Add new row 3
$ws.Range("A3").Value = "sCs"
$ws.Range("B3").Value = "Col9a1"
$ws.Range("C3").Value = "Mag"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.1509656666666667
$ws.Range("H3").Value = 0.452897
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 2
$ws.Range("L3").Value = 0.6666666666666666
$ws.Range("M3").Value = 0.9381496666666668
$ws.Range("N3").Value = 2.814449
$ws.Range("O3").Value = 0.5433136653246862
$ws.Range("P3").Value = 0.5433136653246862
$ws.Range("Q3").Value = 0.1416283898614445
$ws.Range("R3").Value = 1.274655508753
$ws.Range("S3").Value = 0.5433136653246862
$ws.Range("T3").Value = 0.5433136653246862
